$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.2442
$ws.Range("B12").Value = 5.103899999999998
$ws.Range("E14").Value = 16.8758
$ws.Range("E26").Value = 15.9434
$ws.Range("E31").Value = 16.5502
$ws.Range("B32").Value = 6.281299999999999
$ws.Range("E35").Value = 16.77789999999999
$ws.Range("B36").Value = 8.694400000000005
$ws.Range("E37").Value = 16.7725
$ws.Range("B38").Value = 5.094999999999999
$ws.Range("E45").Value = 16.4415
$ws.Range("B46").Value = 5.971500000000004
$ws.Range("B54").Value = 4.599499999999998
$ws.Range("B55").Value = 5.940499999999997
$ws.Range("E57").Value = 16.67020000000001
$ws.Range("B67").Value = 6.090899999999997
$ws.Range("B69").Value = 5.343999999999998
$ws.Range("B72").Value = 5.068200000000004
$ws.Range("B91").Value = 5.750100000000002
$ws.Range("B99").Value = 4.659499999999999
$ws.Range("E100").Value = 16.4083
$ws.Range("E102").Value = 16.62939999999999
